$d = $word.ActiveDocument

# --- Paragraph 1: "Purpose" paragraph ---

# "deduced" -> "deducted"
$d.Content.Find.Execute("deduced from the inventory", $true, $false, $false, $false, $false, $true, 1, $false, "deducted from the inventory", 2) | Out-Null

# "object oriented" -> "object-oriented"
$d.Content.Find.Execute("through object oriented programming", $true, $false, $false, $false, $false, $true, 1, $false, "through object-oriented programming", 2) | Out-Null

# "methods that is used" -> "methods that are used"
$d.Content.Find.Execute("full set of methods that is used", $true, $false, $false, $false, $false, $true, 1, $false, "full set of methods that are used", 2) | Out-Null

# "a database, that stores" -> "a database that stores" (remove comma)
$d.Content.Find.Execute("direct access to a database, that stores", $true, $false, $false, $false, $false, $true, 1, $false, "direct access to a database that stores", 2) | Out-Null

# --- Paragraph 2: "Main interface" paragraph ---

# "employeMenu" -> "employeeMenu" (typo fix)
$d.Content.Find.Execute("employeMenu provides full access", $true, $false, $false, $false, $false, $true, 1, $false, "employeeMenu provides full access", 2) | Out-Null

# "customerMenu provides all access" -> "customerMenu provides access"
$d.Content.Find.Execute("customerMenu provides all access to the functions", $true, $false, $false, $false, $false, $true, 1, $false, "customerMenu provides access to the functions", 2) | Out-Null

# --- Paragraph 3: "Commit Alterations" paragraph ---

# "or "Ad " -> "or "Add " (fix truncated word, split across a page-break run)
$d.Content.Find.Execute("or “Ad ", $true, $false, $false, $false, $false, $true, 1, $false, "or “Add ", 2) | Out-Null
